$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Red Status: 7 projects"

$ws.Range("G4").Value = "c1: 15"
$ws.Range("G5").Value = "c2: 16"
$ws.Range("G6").Value = "c3: 17"
$ws.Range("G7").Value = "o1: 18"
$ws.Range("G8").Value = "o2: 19"
